$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.8.2"
$wsMeta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

$wsElem = $wb.Worksheets.Item("Elements")
$constraint = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
$wsElem.Range("AJ1").Value = $constraint
$wsElem.Range("AJ3").Value = $constraint
